$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark. In the original document it
#     sits mid-paragraph, right after "...triggered dialogue.   " (inside the
#     "The second aspect of audio is dialogue..." paragraph). ---
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# --- Step 2: rewrite the "Objective" paragraph. The original paragraph had
#     two runs:
#       "The main objective of the project is to create a scene/level to
#        showcase the elements discussed above" + ", I'll explore and
#        demonstrate these techniques in depth."
#     -> becomes a single run with the new wording. ---
$rng = $d.Content
$rng.Find.ClearFormatting()
$oldText = "The main objective of the project is to create a scene/level to showcase the elements discussed above, I" + [char]8217 + "ll explore and demonstrate these techniques in depth."
$newText = "The overall objective of this project is too create a UE4 project that showcases the above elements in depth, too increase my overall knowledge of Unreal games development."
$rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# --- Step 3: add the "_GoBack" bookmark back, now collapsed right at the end
#     of that same paragraph's text (after the new run, before the paragraph
#     mark). $rng, after a successful Find/Replace, already refers to the
#     freshly-inserted replacement text, so we can anchor off it directly.
#     A bare Bookmarks.Add() at the "end of paragraph" position collapses
#     oddly in this host, so insert a throwaway marker char right after the
#     target spot, bookmark just before the marker, then delete the marker
#     -- the bookmark stays put exactly where we wanted it. ---
$endRng = $rng.Duplicate
$endRng.Collapse(0) | Out-Null      # wdCollapseEnd -- right after the new text

$endRng.InsertAfter("X")
$markerPos = $endRng.Start
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()
